$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-19 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("718÷2=359, 0", $true, $false, $false, $false, $false, $true, 1, $false, "945÷8=118, 1", 2) | Out-Null
$d.Content.Find.Execute("350÷7=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "844÷4=211, 0", 2) | Out-Null
$d.Content.Find.Execute("394÷8=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "705÷9=78, 3", 2) | Out-Null
$d.Content.Find.Execute("219÷4=54, 3", $true, $false, $false, $false, $false, $true, 1, $false, "954÷5=190, 4", 2) | Out-Null
$d.Content.Find.Execute("386÷6=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "278÷8=34, 6", 2) | Out-Null
$d.Content.Find.Execute("998÷8=124, 6", $true, $false, $false, $false, $false, $true, 1, $false, "252÷8=31, 4", 2) | Out-Null
$d.Content.Find.Execute("616÷8=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "366÷4=91, 2", 2) | Out-Null
$d.Content.Find.Execute("829÷3=276, 1", $true, $false, $false, $false, $false, $true, 1, $false, "322÷5=64, 2", 2) | Out-Null
$d.Content.Find.Execute("659÷5=131, 4", $true, $false, $false, $false, $false, $true, 1, $false, "934÷8=116, 6", 2) | Out-Null
$d.Content.Find.Execute("661÷7=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "800÷9=88, 8", 2) | Out-Null
$d.Content.Find.Execute("727÷7=103, 6", $true, $false, $false, $false, $false, $true, 1, $false, "556÷2=278, 0", 2) | Out-Null
$d.Content.Find.Execute("569÷4=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "164÷4=41, 0", 2) | Out-Null
$d.Content.Find.Execute("110÷5=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "594÷9=66, 0", 2) | Out-Null
$d.Content.Find.Execute("539÷2=269, 1", $true, $false, $false, $false, $false, $true, 1, $false, "748÷4=187, 0", 2) | Out-Null
$d.Content.Find.Execute("102÷5=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "317÷5=63, 2", 2) | Out-Null
$d.Content.Find.Execute("740÷7=105, 5", $true, $false, $false, $false, $false, $true, 1, $false, "666÷5=133, 1", 2) | Out-Null
$d.Content.Find.Execute("943÷6=157, 1", $true, $false, $false, $false, $false, $true, 1, $false, "413÷5=82, 3", 2) | Out-Null
$d.Content.Find.Execute("699÷3=233, 0", $true, $false, $false, $false, $false, $true, 1, $false, "961÷7=137, 2", 2) | Out-Null
$d.Content.Find.Execute("573÷3=191, 0", $true, $false, $false, $false, $false, $true, 1, $false, "994÷2=497, 0", 2) | Out-Null
$d.Content.Find.Execute("523÷6=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷3=84, 1", 2) | Out-Null
$d.Content.Find.Execute("117÷2=58, 1", $true, $false, $false, $false, $false, $true, 1, $false, "463÷9=51, 4", 2) | Out-Null
$d.Content.Find.Execute("392÷7=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "470÷2=235, 0", 2) | Out-Null
$d.Content.Find.Execute("941÷4=235, 1", $true, $false, $false, $false, $false, $true, 1, $false, "426÷3=142, 0", 2) | Out-Null
$d.Content.Find.Execute("941÷3=313, 2", $true, $false, $false, $false, $false, $true, 1, $false, "373÷3=124, 1", 2) | Out-Null
$d.Content.Find.Execute("571÷4=142, 3", $true, $false, $false, $false, $false, $true, 1, $false, "847÷6=141, 1", 2) | Out-Null
